## IEAGHG scenario update: rechecking of emission and energy data based on source
$wb = $excel.ActiveWorkbook

$wsConnections = $wb.Worksheets.Item("IEAGHG connections")
$wsChains      = $wb.Worksheets.Item("IEAGHG chains")
$wsSteelChain  = $wb.Worksheets.Item("IEAGHG steel chain")

# ---------------------------------------------------------------------------
# "IEAGHG connections" sheet updates
# ---------------------------------------------------------------------------

# Relabel the blast-furnace-gas energy product name (row 10, "o product" column)
$wsConnections.Range("E10").Value = "energy in blast furnace gas"
$wsConnections.Range("E10").Style = "Normal"

# Add the newly-rechecked purge % value for the last connection row
$wsConnections.Range("K19").Value = 0.68

# New (blank, but formatted) row underneath the table - mirrors the formats
# used on the row above it
$wsConnections.Range("C19").Copy() | Out-Null
$wsConnections.Range("C20").PasteSpecial(-4122) | Out-Null

$wsConnections.Range("G19").Copy() | Out-Null
$wsConnections.Range("E20").PasteSpecial(-4122) | Out-Null
$wsConnections.Range("F20").PasteSpecial(-4122) | Out-Null
$wsConnections.Range("G20").PasteSpecial(-4122) | Out-Null
$wsConnections.Range("H20").PasteSpecial(-4122) | Out-Null
$wsConnections.Range("I20").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# "IEAGHG chains" sheet updates
# ---------------------------------------------------------------------------

# New chain: flare -> heat_flare (outflow / inflow mapping for the new
# "flare" energy loss accounted for in the rechecked balance)
$wsChains.Range("A8").Value = "flare"
$wsChains.Range("B8").Value = "heat"
$wsChains.Range("C8").Value = "inflow"
$wsChains.Range("E8").Value = "heat_flare"
$wsChains.Range("E8").NumberFormat = "@"

# ---------------------------------------------------------------------------
# Restore the view state Excel saves on exit (active sheet / selections)
# ---------------------------------------------------------------------------

$wsChains.Range("E8").Select() | Out-Null

$wsConnections.Range("M1").Select() | Out-Null
$wsConnections.Activate() | Out-Null

$wsSteelChain.Range("C11").Select() | Out-Null

Write-Host "IEAGHG factories workbook updated"
